# Update the multiplication answers in the worksheet table.
$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Unique old -> new pairs (safe to use a document-wide Find & Replace).
Replace-Exact "206×7=1442" "735×9=6615"
Replace-Exact "747×2=1494" "987×4=3948"
Replace-Exact "113×7=791" "435×5=2175"
Replace-Exact "549×3=1647" "156×4=624"
Replace-Exact "446×5=2230" "408×8=3264"
Replace-Exact "183×8=1464" "289×8=2312"
Replace-Exact "838×5=4190" "361×8=2888"
Replace-Exact "946×3=2838" "936×4=3744"
Replace-Exact "467×2=934" "750×3=2250"
Replace-Exact "668×6=4008" "465×2=930"
Replace-Exact "987×6=5922" "342×4=1368"
Replace-Exact "957×7=6699" "385×4=1540"
Replace-Exact "814×6=4884" "840×6=5040"
Replace-Exact "456×4=1824" "460×2=920"
Replace-Exact "574×9=5166" "204×2=408"
Replace-Exact "427×7=2989" "279×4=1116"
Replace-Exact "364×3=1092" "538×3=1614"
Replace-Exact "142×6=852" "422×7=2954"
Replace-Exact "457×9=4113" "754×6=4524"
Replace-Exact "129×8=1032" "160×8=1280"
Replace-Exact "573×2=1146" "802×7=5614"
Replace-Exact "909×4=3636" "652×9=5868"
Replace-Exact "559×8=4472" "675×9=6075"

# "881×7=6167" occurs twice and must become two different values, so
# target each table cell directly instead of using a global replace.
$t = $d.Tables(1)
$t.Cell(15, 2).Range.Text = "121×9=1089"
$t.Cell(20, 4).Range.Text = "899×2=1798"
